# Update message in controller
# Append newly redeemed chest codes to the "Codes" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codes")

$newCodes = @(
    "shak-adha-nidb-2022",
    "OLEO-PEAK-COWY",
    "BLUE-CARB-ZILL",
    "KOBO-LDSY-AY!!",
    "SOPS-HORE-HAPS",
    "TUNE-INTM-RROW",
    "WALN-UTIN-SPAA-ACE!",
    "MARSINOFFICE",
    "EXTR-ALIF-EEVE ",
    "SPAG-PENS-DEVA",
    "SKETCHEXTRA!",
    "HOPS-APSE-FOLK",
    "FANO-FVIR-GIL!",
    "LOOK-MORE-FORM",
    "MULL-CLEW-SALE",
    "YODH-GOWL-LEST",
    "ICES-GAMY-PIKI",
    "ALLS-FRAY-SPIF-WALL",
    "FEASTONFORMS",
    "ULEX-RORT-MASU",
    "ONYO-UTUB-EEP5",
    "XTRALIFESOON",
    "JCMR-AIDT-IME!"
)

# Find the first empty row right after the existing data (row 864 given 863 rows already used).
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newCodes.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newCodes[$i]
    $ws.Cells.Item($row, 2).Value = $true
}

# Restore the view to show the newly added rows, mirroring the author's scroll/selection.
$excel.ActiveWindow.ScrollRow = 830
$excel.ActiveWindow.ScrollColumn = 1
$ws.Cells.Item($startRow, 2).Select()
